# Update the workbook: set cell A1 on the active sheet and configure the
# page setup (paper size / orientation) to match the authored change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "543153ss"

$ps = $ws.PageSetup
$ps.PaperSize = 9      # xlPaperA4
$ps.Orientation = 1    # xlPortrait
